# Added Mockup Task, because it's almost done
# Set the assignee initials ("TA") for the "Mockup creation" task row (row 13).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")
$ws.Range("C13").Value = "TA"
